$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Move the "3" story-point value from D13 to E13 (it shifted one day later)
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 3

# Move the "2" story-point value from J13 to L13 (it shifted two days later)
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 2

# Update the view state to match where the user was working
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("E15").Select()

$wb.Save()
